# Appends new sensor-log rows to the SeniorConnect master log workbook,
# matching the "2026-01-30 18:39 - 18:41" batch of new readings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data note: every new row's Date column is the literal text
# "2026-01-30" and every new row's Hour column is the literal text
# "18:00". Columns that look like dates/times/percentages must have
# their NumberFormat forced to Text ("@") before the value is written,
# otherwise Excel will silently re-interpret "2026-01-30" as a date
# serial, "85.0%" as a fraction, etc.
# ---------------------------------------------------------------------

# ========================= ALERTS (sheet1) ============================
$ws = $wb.Worksheets.Item("ALERTS")

$ws.Range("A12:A13").NumberFormat = "@"
$ws.Range("A12:A13").Value = "2026-01-30"
$ws.Range("C12:C13").Value = "18:00"
$ws.Range("D12:D13").Value = "Living Room"
$ws.Range("E12:E13").Value = "CRITICAL"
$ws.Range("F12:F13").Value = "FALL_DETECTED"

$ws.Cells.Item(12, 2).Value = "18:40:10"
$ws.Cells.Item(13, 2).Value = "18:40:14"

# ========================= PIR (sheet2) ================================
$ws = $wb.Worksheets.Item("PIR")

$firstRow = 221
$lastRow = 237

$ws.Range("A" + $firstRow + ":A" + $lastRow).NumberFormat = "@"
$ws.Range("A" + $firstRow + ":A" + $lastRow).Value = "2026-01-30"
$ws.Range("C" + $firstRow + ":C" + $lastRow).Value = "18:00"
$ws.Range("D" + $firstRow + ":D" + $lastRow).Value = "Bathroom"
$ws.Range("E" + $firstRow + ":E" + $lastRow).Value = "No Motion"
$ws.Range("F" + $firstRow + ":F" + $lastRow).Value = "Inactive"

$pirTimes = @(
    "18:39:36","18:39:38","18:39:43","18:39:48","18:39:53",
    "18:40:14","18:40:17","18:40:19","18:40:23","18:40:28",
    "18:40:33","18:40:38","18:40:43","18:40:48","18:40:53",
    "18:40:58","18:41:03"
)
for ($i = 0; $i -lt $pirTimes.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 2).Value = $pirTimes[$i]
}

# ========================= Humidity (sheet3) ===========================
$ws = $wb.Worksheets.Item("Humidity")

$firstRow = 145
$lastRow = 154

$ws.Range("A" + $firstRow + ":A" + $lastRow).NumberFormat = "@"
$ws.Range("A" + $firstRow + ":A" + $lastRow).Value = "2026-01-30"
$ws.Range("C" + $firstRow + ":C" + $lastRow).Value = "18:00"
$ws.Range("D" + $firstRow + ":D" + $lastRow).Value = "Bathroom"
$ws.Range("E" + $firstRow + ":E" + $lastRow).NumberFormat = "@"
$ws.Range("F" + $firstRow + ":F" + $lastRow).Value = "Active"

$humTimes = @(
    "18:39:37","18:39:38","18:39:48","18:39:53","18:40:15",
    "18:40:19","18:40:28","18:40:33","18:40:49","18:40:59"
)
$humValues = @(
    "85.0%","86.0%","86.0%","85.1%","86.1%",
    "86.1%","86.1%","85.1%","86.1%","86.0%"
)
for ($i = 0; $i -lt $humTimes.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 2).Value = $humTimes[$i]
    $ws.Cells.Item($firstRow + $i, 5).Value = $humValues[$i]
}

# ========================= Proximity (sheet5) ===========================
$ws = $wb.Worksheets.Item("Proximity")

$firstRow = 15
$lastRow = 21

$ws.Range("A" + $firstRow + ":A" + $lastRow).NumberFormat = "@"
$ws.Range("A" + $firstRow + ":A" + $lastRow).Value = "2026-01-30"
$ws.Range("C" + $firstRow + ":C" + $lastRow).Value = "18:00"
$ws.Range("D" + $firstRow + ":D" + $lastRow).Value = "Living Room Main Door"

$proxTimes = @("18:40:15","18:40:18","18:40:29","18:40:40","18:40:44","18:40:55","18:41:01")
$proxKinds = @("ENTER","EXIT","ENTER","EXIT","ENTER","EXIT","ENTER")
for ($i = 0; $i -lt $proxTimes.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 2).Value = $proxTimes[$i]
    $ws.Cells.Item($row, 5).Value = $proxKinds[$i]
    if ($proxKinds[$i] -eq "ENTER") {
        $ws.Cells.Item($row, 6).Value = "User ENTERED Living Room Main Door"
    } else {
        $ws.Cells.Item($row, 6).Value = "User EXITED Living Room Main Door"
    }
}

# ========================= Camera (sheet7) ===============================
$ws = $wb.Worksheets.Item("Camera")

$firstRow = 15
$lastRow = 27

$ws.Range("A" + $firstRow + ":A" + $lastRow).NumberFormat = "@"
$ws.Range("A" + $firstRow + ":A" + $lastRow).Value = "2026-01-30"
$ws.Range("C" + $firstRow + ":C" + $lastRow).Value = "18:00"
$ws.Range("D" + $firstRow + ":D" + $lastRow).Value = "Living Room Main Door"
$ws.Range("E" + $firstRow + ":E" + $lastRow).Value = "Image Captured"
$ws.Range("F" + $firstRow + ":F" + $lastRow).Value = "Active"

$camTimes = @(
    "18:40:15","18:40:16","18:40:16","18:40:17","18:40:18",
    "18:40:18","18:40:20","18:40:30","18:40:42","18:40:45",
    "18:40:58","18:41:02","18:41:05"
)
for ($i = 0; $i -lt $camTimes.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 2).Value = $camTimes[$i]
}
